$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BESS")

# Clear the user-input values in column C (rows 2-38), leaving formulas/labels intact.
$ws.Range("C2:C38").ClearContents()

# Rows 5, 18 and 33 had their custom (auto-fit) heights driven by now-removed
# wrapped text; re-run AutoFit so they collapse back to the default height.
$ws.Rows("5:5").AutoFit()
$ws.Rows("18:18").AutoFit()
$ws.Rows("33:33").AutoFit()

# Move the active selection to C7 (matches the post-edit selection in the file).
$ws.Range("C7").Select()
